$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 36: Monster Hunter Portable 2nd
$ws.Range("A36").Value = "Monster Hunter Portable 2nd"
$ws.Range("B36").Value = "JPN"
$ws.Range("C36").Value = 39135
$ws.Range("C36").NumberFormat = "d-mmm-yy"
$ws.Range("D36").Value = 2
$ws.Range("E36").Value = "Yasunori Ichinose"
$ws.Range("F36").Value = "Playstation Portable"

# Row 37: Monster Hunter Portable 2nd G
$ws.Range("A37").Value = "Monster Hunter Portable 2nd G"
$ws.Range("B37").Value = "JPN"
$ws.Range("C37").Value = 39534
$ws.Range("C37").NumberFormat = "d-mmm-yy"
$ws.Range("D37").Value = 2
$ws.Range("E37").Value = "Yasunori Ichinose"
$ws.Range("F37").Value = "Playstation Portable"

# Row 38: Monster Hunter 3 G
$ws.Range("A38").Value = "Monster Hunter 3 G"
$ws.Range("B38").Value = "JPN"
$ws.Range("C38").Value = 40887
$ws.Range("C38").NumberFormat = "d-mmm-yy"
$ws.Range("D38").Value = 3
$ws.Range("E38").Value = "Kaname Fujioka"
$ws.Range("F38").Value = "3DS"

# Update selection to match final state
$ws.Range("A18").Select()
